$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 309, shifting existing rows (309..341) down to (310..342).
$ws.Rows("309:309").Insert()

# Populate the newly inserted row 309 with a new weekly record (same market/product
# metadata as the row that used to sit at 309, but a new date + new price figures).
$ws.Range("A309").Value = 10
$ws.Range("B309").Value = "Vega Modelo de Temuco"
$ws.Range("C309").Value = "La Araucanía"
$ws.Range("D309").Value = 45223
$ws.Range("E309").Value = 9
$ws.Range("F309").Value = "Fruta"
$ws.Range("G309").Value = 100101
$ws.Range("H309").Value = "Berries"
$ws.Range("I309").Value = 100112025
$ws.Range("J309").Value = "Frutilla"
$ws.Range("K309").Value = "Sin especificar"
$ws.Range("L309").Value = "Primera"
$ws.Range("M309").Value = 35
$ws.Range("N309").Value = 12000
$ws.Range("O309").Value = 12000
$ws.Range("P309").Value = 12000
$ws.Range("Q309").Value = "$/bandeja 7 kilos"
$ws.Range("R309").Value = "Provincia de Melipilla"
$ws.Range("S309").Value = 1714
$ws.Range("T309").Value = 7
